# Change "delete" to "kill" in sequence diagram
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1. TextBox 23: "delete 1" -> "kill 1"
$deleteShape = $s.Shapes.Item("TextBox 23")
$deleteShape.TextFrame.TextRange.Text = "kill 1"

# 2. TextBox 28: "killTasks(p)" - the trailing "(p" and ")" runs get merged
#    into a single "(p)" run (same formatting), matching the canonical edit.
$killTasksShape = $s.Shapes.Item("TextBox 28")
$tr = $killTasksShape.TextFrame.TextRange
$sub = $tr.Characters(10, 3)
$sub.Text = "(p)"
